$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.156.64'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '1.842.30'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6881'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3017'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07467'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.21%  '
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07659'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("D12").Value = '1.844.42'
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.066'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6838'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.62'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.181'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.42%  '
$ws.Range("D17").Value = '29.158.33'
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008171'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.72%  '
$ws.Range("D19").Value = '2.076.47'
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '228.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.400'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9996'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("E25").Value = '  -3.70%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.782'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.515'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.145'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.195'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05257'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7596'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.856'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.39%  '
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("D38").Value = '1.305.97'
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01833'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.724'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9312'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.942'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.97'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9992'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("E45").Value = '  +0.90%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000123'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5195'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.977.41'
$ws.Range("E48").Value = '  -0.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.509'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.774'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05961'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.02%  '
